$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    2 = 2500
    3 = 1000
    4 = 750
    11 = 2000
    13 = 1000
    14 = 2500
    18 = 3000
    19 = 750
    21 = 1000
    22 = 1500
    26 = 4000
    28 = 750
    29 = 500
    31 = 1000
    32 = 0
    35 = 8200
    36 = 3000
    41 = 3500
    42 = 2500
    49 = 3500
    51 = 3000
    52 = 3500
    54 = 3500
    56 = 2000
    57 = 8000
    58 = 3500
    59 = 3500
    60 = 1000
    61 = 3500
    64 = 3500
    65 = 3500
    66 = 5501
    67 = 5800
    69 = 4000
    71 = 3000
    73 = 5800
    74 = 4000
    75 = 3000
    77 = 2300
    78 = 2500
    79 = 3200
    80 = 4300
    82 = 500
    83 = 2100
    91 = 1000
    92 = 15000
    93 = 10000
    94 = 1500
    96 = 500
    97 = 200
    98 = 1500
    99 = 1000
    102 = 3000
    103 = 6000
    104 = 6000
    106 = 3000
    107 = 1000
    108 = 500
    109 = 1000
    116 = 1500
    117 = 6000
    119 = 1500
    120 = 3000
    121 = 300
    123 = 10000
    124 = 8000
    125 = 200
    128 = 3500
    129 = 1000
    130 = 2000
    131 = 200
    133 = 300
    134 = 1500
    135 = 1500
    136 = 1500
    137 = 500
    140 = 1500
    141 = 500
    142 = 1000
    144 = 500
    146 = 4000
    147 = 4000
    148 = 1500
    151 = 6000
    153 = 5000
    154 = 6500
    155 = 500
    156 = 500
    157 = 8000
    159 = 4500
    162 = 5000
    163 = 2000
    164 = 3500
    165 = 2000
    167 = 2000
    170 = 3000
    171 = 10000
    173 = 3500
    175 = 1000
    177 = 2000
    180 = 500
    183 = 1000
    187 = 3000
    188 = 2000
    199 = 0
    200 = 1000
    205 = 1000
    206 = 4000
    210 = 0
    212 = 2000
    214 = 1000
    216 = 0
    220 = 1000
    222 = 2000
    223 = 0
    226 = 0
    229 = 0
    235 = 6000
    236 = 2000
    237 = 0
    238 = 2500
    243 = 2000
    255 = 500
    256 = 3500
    257 = 250
    258 = 2000
    259 = 12000
    261 = 2000
    262 = 1425
    263 = 12000
    264 = 3000
    265 = 250
    267 = 1500
    268 = 1500
    270 = 300
    271 = 20000
    281 = 440476
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $updates[$row]
}

Write-Host "Done updating PRESUPUESTO column (G) on VENTA MENSUAL sheet."
